$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 = "Save" - copy the formatting from the existing
# header cell G1 (bold, bordered, centered) so no new style is minted.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data cell H2 = 0
$ws.Range("H2").Value = 0
